$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.082199919036966662
$ws.Range("A2").Value = -0.009999999617175348
$ws.Range("A3").Value = -0.0089999996210288202
$ws.Range("A4").Value = 0.28399342830664409
$ws.Range("A5").Value = -0.0059999996314914483
$ws.Range("A6").Value = -0.0059999996172948045
$ws.Range("A7").Value = -0.019999999547087199
$ws.Range("A8").Value = -0.019999999543006908
$ws.Range("A9").Value = -0.0059999996071757877
$ws.Range("A10").Value = -0.005999999603155004
$ws.Range("A11").Value = -0.0044999996104451156
$ws.Range("A12").Value = -0.005580344626592737
$ws.Range("A13").Value = -0.013488974724367608
$ws.Range("A14").Value = -0.011999999568490516
$ws.Range("A15").Value = -0.0059999995978685661
$ws.Range("A16").Value = 0.019240816393325755
$ws.Range("A17").Value = -0.0059999995964004071
$ws.Range("A18").Value = -0.0089999995812855005
$ws.Range("A19").Value = -0.0089999996214200628
$ws.Range("A20").Value = -0.0089999996156837625
$ws.Range("A21").Value = -0.026784856280380076
$ws.Range("A22").Value = -0.0089999996142031691
$ws.Range("A23").Value = -0.0089999996176919339
$ws.Range("A24").Value = -0.041999999449444658
$ws.Range("A25").Value = -0.041999999446503011
$ws.Range("A26").Value = -0.0059999996159412206
$ws.Range("A27").Value = -0.0059999996140902567
$ws.Range("A28").Value = -0.0059999996064119543
$ws.Range("A29").Value = -0.011999999571489894
$ws.Range("A30").Value = -0.019999999529681567
$ws.Range("A31").Value = -0.013647289696809395
$ws.Range("A32").Value = -0.020999999520372903
$ws.Range("A33").Value = -0.0059999995938486705
